# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N on the "Repayment schedule"
# sheet -- shifting the old N/O/P columns (Late / heading / Outstanding)
# one column to the right -- and makes "Repayment schedule" the active
# sheet/tab instead of "Transactions".

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (existing N, O, P shift right to O, P, Q).
$wsSchedule.Range("N1").EntireColumn.Insert()

# The newly inserted column ends up a bit wider than its neighbours by
# default; set it to display as width 11 (same look as column M).
$wsSchedule.Columns("N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab, with K20 selected,
# instead of "Transactions" (which loses the tab-selected / active state).
$wsSchedule.Activate()
$wsSchedule.Range("K20").Select()
